$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.255.43"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").Value = "1.596.95"
$ws.Range("E3").Value = "  +0.56%  "

$ws.Range("D5").Value = "'211.47"
$ws.Range("E5").Value = "  -0.16%  "

$ws.Range("D6").Value = "'0.505"
$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  +0.19%  "

$ws.Range("E10").Value = "  -1.29%  "

$ws.Range("E11").Value = "  +0.71%  "

$ws.Range("D12").Value = "1.821.81"
$ws.Range("E12").Value = "  +0.61%  "

$ws.Range("D13").Value = "1.614.95"
$ws.Range("E13").Value = "  +1.55%  "

$ws.Range("E14").Value = "  -0.40%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.503"
$ws.Range("E15").Value = "  -2.42%  "

$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "'63.67"
$ws.Range("E16").Value = "  -0.39%  "

$ws.Range("D17").Value = "26.264.75"
$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("D18").Value = "'230.33"
$ws.Range("E18").Value = "  +7.55%  "

$ws.Range("D19").Value = "'7.67"
$ws.Range("E19").Value = "  +4.01%  "

$ws.Range("E20").Value = "  -0.51%  "

$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("D22").Value = "'4.22"
$ws.Range("E22").Value = "  -0.55%  "

$ws.Range("D23").Value = "'8.94"
$ws.Range("E23").Value = "  -0.41%  "

$ws.Range("E24").Value = "  +1.14%  "

$ws.Range("D25").Value = "'146.15"
$ws.Range("E25").Value = "  +0.99%  "

$ws.Range("D28").Value = "'0.113"
$ws.Range("E28").Value = "  +0.32%  "

$ws.Range("E29").Value = "  +1.39%  "

$ws.Range("E30").Value = "  -0.51%  "

$ws.Range("E31").Value = "  -0.08%  "

$ws.Range("E32").Value = "  +0.59%  "

$ws.Range("D33").Value = "1.466.18"
$ws.Range("E33").Value = "  +2.96%  "

$ws.Range("E34").Value = "  +0.07%  "

$ws.Range("E35").Value = "  -0.35%  "

$ws.Range("D37").Value = "'0.568"
$ws.Range("E37").Value = "  -3.41%  "

$ws.Range("E38").Value = "  -1.10%  "

$ws.Range("D39").Value = "'0.821"
$ws.Range("E39").Value = "  -0.20%  "

$ws.Range("D40").Value = "'5.76"
$ws.Range("E40").Value = "  -2.57%  "

$ws.Range("E42").Value = "  +2.38%  "

$ws.Range("E43").Value = "  -1.57%  "

$ws.Range("D44").Value = "1.734.40"
$ws.Range("E44").Value = "  +0.72%  "

$ws.Range("D45").Value = "'0.756"
$ws.Range("E45").Value = "  -1.32%  "

$ws.Range("D46").Value = "'60.51"
$ws.Range("E46").Value = "  -1.02%  "

$ws.Range("D47").Value = "'87.87"
$ws.Range("E47").Value = "  +2.47%  "

$ws.Range("E48").Value = "  -0.83%  "

$ws.Range("E49").Value = "  -0.10%  "

$ws.Range("D50").Value = "'0.998"
$ws.Range("E50").Value = "  +0.01%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.41"
$ws.Range("E51").Value = "  +0.96%  "
